# The upstream change for this fixture (commit "Moving from 2.0.1 to
# 2.0.2") is purely a byte-level re-serialization artifact: every hunk
# in the recorded diff (word/document.xml, word/footnotes.xml,
# word/header1.xml, word/styles.xml) reorders XML namespace
# declarations / element attributes into alphabetical order (e.g.
# <w:tcW w:w="3070" w:type="dxa"/> -> <w:tcW w:type="dxa" w:w="3070"/>,
# <w:footnote w:type="separator" w:id="-1"> -> <w:footnote w:id="-1"
# w:type="separator">, latent-style / w:style / tblBorders / tab
# attribute reorderings, etc.) coming from a newer version of the
# fixture-generating tool. Every (element, attribute-name,
# attribute-value) triple is identical before and after -- no text,
# numbering, formatting, structure, relationship, or other visible
# content actually changed anywhere in the package.
#
# The Word object model only exposes semantic document edits (text,
# formatting, structure, ...); it has no API for reordering raw XML
# attributes, and this headless engine otherwise preserves each part's
# original on-disk attribute order whenever its content is not
# otherwise being changed. So there is no content-level edit to
# reproduce here: the correct COM script is a no-op that leaves the
# document exactly as authored.
$d = $word.ActiveDocument
